$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 45.76217133333333
$ws.Range("H2").Value = 137.286514
$ws.Range("I2").Value = 0.6763939203605134
$ws.Range("J2").Value = 0.6763939203605135
$ws.Range("M2").Value = 2.718682666666667
$ws.Range("N2").Value = 8.156048
$ws.Range("O2").Value = 0.5434637507613679
$ws.Range("P2").Value = 0.5434637507613679
$ws.Range("Q2").Value = 124.4128219929636
$ws.Range("R2").Value = 1119.715397936672
$ws.Range("S2").Value = 0.3675955769513106
$ws.Range("T2").Value = 0.3675955769513106

$ws.Range("G3").Value = 45.76217133333333
$ws.Range("H3").Value = 137.286514
$ws.Range("I3").Value = 0.6763939203605134
$ws.Range("J3").Value = 0.6763939203605135
$ws.Range("M3").Value = 2.283827
$ws.Range("N3").Value = 6.851481
$ws.Range("O3").Value = 0.4565362492386322
$ws.Range("P3").Value = 0.4565362492386321
$ws.Range("Q3").Value = 104.5128824696927
$ws.Range("R3").Value = 940.6159422272341
$ws.Range("S3").Value = 0.3087983434092029
$ws.Range("T3").Value = 0.3087983434092029

$ws.Range("I4").Value = 0.1388778842960613
$ws.Range("J4").Value = 0.1388778842960613
$ws.Range("M4").Value = 2.718682666666667
$ws.Range("N4").Value = 8.156048
$ws.Range("O4").Value = 0.5434637507613679
$ws.Range("P4").Value = 0.5434637507613679
$ws.Range("Q4").Value = 25.54456652785422
$ws.Range("R4").Value = 229.901098750688
$ws.Range("S4").Value = 0.07547509589734074
$ws.Range("T4").Value = 0.07547509589734076

$ws.Range("I5").Value = 0.1388778842960613
$ws.Range("J5").Value = 0.1388778842960613
$ws.Range("M5").Value = 2.283827
$ws.Range("N5").Value = 6.851481
$ws.Range("O5").Value = 0.4565362492386322
$ws.Range("P5").Value = 0.4565362492386321
$ws.Range("Q5").Value = 21.45869080452066
$ws.Range("R5").Value = 193.128217240686
$ws.Range("S5").Value = 0.06340278839872056
$ws.Range("T5").Value = 0.06340278839872056

$ws.Range("G6").Value = 12.29750866666667
$ws.Range("H6").Value = 36.892526
$ws.Range("I6").Value = 0.1817649787009828
$ws.Range("J6").Value = 0.1817649787009828
$ws.Range("M6").Value = 2.718682666666667
$ws.Range("N6").Value = 8.156048
$ws.Range("O6").Value = 0.5434637507613679
$ws.Range("P6").Value = 0.5434637507613679
$ws.Range("Q6").Value = 33.43302365524978
$ws.Range("R6").Value = 300.897212897248
$ws.Range("S6").Value = 0.09878267708189624
$ws.Range("T6").Value = 0.09878267708189624

$ws.Range("G7").Value = 12.29750866666667
$ws.Range("H7").Value = 36.892526
$ws.Range("I7").Value = 0.1817649787009828
$ws.Range("J7").Value = 0.1817649787009828
$ws.Range("M7").Value = 2.283827
$ws.Range("N7").Value = 6.851481
$ws.Range("O7").Value = 0.4565362492386322
$ws.Range("P7").Value = 0.4565362492386321
$ws.Range("Q7").Value = 28.08538232566734
$ws.Range("R7").Value = 252.768440931006
$ws.Range("S7").Value = 0.08298230161908654
$ws.Range("T7").Value = 0.08298230161908653

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2004796666666666
$ws.Range("H8").Value = 0.601439
$ws.Range("I8").Value = 0.002963216642442438
$ws.Range("J8").Value = 0.002963216642442439
$ws.Range("M8").Value = 2.718682666666667
$ws.Range("N8").Value = 8.156048
$ws.Range("O8").Value = 0.5434637507613679
$ws.Range("P8").Value = 0.5434637507613679
$ws.Range("Q8").Value = 0.5450405947857777
$ws.Range("R8").Value = 4.905365353072
$ws.Range("S8").Value = 0.001610400830820274
$ws.Range("T8").Value = 0.001610400830820275

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2004796666666666
$ws.Range("H9").Value = 0.601439
$ws.Range("I9").Value = 0.002963216642442438
$ws.Range("J9").Value = 0.002963216642442439
$ws.Range("M9").Value = 2.283827
$ws.Range("N9").Value = 6.851481
$ws.Range("O9").Value = 0.4565362492386322
$ws.Range("P9").Value = 0.4565362492386321
$ws.Range("Q9").Value = 0.4578608756843333
$ws.Range("R9").Value = 4.120747881159
$ws.Range("S9").Value = 0.001352815811622164
$ws.Range("T9").Value = 0.001352815811622164
